$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)
$p.Range.Text = "2024-05-11 Saturday"

$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "31÷6="
$t.Cell(1,2).Range.Text = "83÷4="
$t.Cell(1,3).Range.Text = "65÷3="
$t.Cell(1,4).Range.Text = "72÷5="
$t.Cell(1,5).Range.Text = "13÷8="
$t.Cell(5,1).Range.Text = "61÷7="
$t.Cell(5,2).Range.Text = "79÷4="
$t.Cell(5,3).Range.Text = "83÷5="
$t.Cell(5,4).Range.Text = "46÷7="
$t.Cell(5,5).Range.Text = "28÷5="
$t.Cell(9,1).Range.Text = "66÷4="
$t.Cell(9,2).Range.Text = "92÷3="
$t.Cell(9,3).Range.Text = "58÷8="
$t.Cell(9,4).Range.Text = "73÷4="
$t.Cell(9,5).Range.Text = "83÷9="
$t.Cell(13,1).Range.Text = "31÷3="
$t.Cell(13,2).Range.Text = "44÷7="
$t.Cell(13,3).Range.Text = "64÷5="
$t.Cell(13,4).Range.Text = "57÷6="
$t.Cell(13,5).Range.Text = "15÷2="
$t.Cell(17,1).Range.Text = "66÷3="
$t.Cell(17,2).Range.Text = "16÷9="
$t.Cell(17,3).Range.Text = "87÷3="
$t.Cell(17,4).Range.Text = "21÷8="
$t.Cell(17,5).Range.Text = "40÷3="
